$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "Celebratory "

# Rows 2-16: not yet celebrated at that point in the workbook's original scheme -> "FALSE " (text)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = "FALSE "
}

# Rows 17-46: mark as celebratory (boolean TRUE)
for ($r = 17; $r -le 46; $r++) {
    $ws.Cells.Item($r, 5).Value = $true
}

# Restore the view: scroll down a bit and select F45, as in the edited workbook
$ws.Activate()
$ws.Range("F45").Select()
